# Scheduled-runner market data refresh for Phoenix_Profits sheets.
# Updates currentAveragePrice*/LevePrice*/LeveProfit* columns (H-N) per leve row
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("H64").Value = 2833.3333
$ws.Range("J64").Value = 3500
$ws.Range("L64").Value = 3500
$ws.Range("N64").Value = -3996

# Row 67
$ws.Range("H67").Value = 2833.3333
$ws.Range("J67").Value = 3500
$ws.Range("L67").Value = 3500
$ws.Range("N67").Value = -5216

# Row 74
$ws.Range("H74").Value = 11466.8
$ws.Range("I74").Value = 6167.1665
$ws.Range("K74").Value = 6167.1665
$ws.Range("M74").Value = -5231.1665

# Row 77
$ws.Range("H77").Value = 11466.8
$ws.Range("I77").Value = 6167.1665
$ws.Range("K77").Value = 30835.8325
$ws.Range("M77").Value = -26155.8325

# Row 132
$ws.Range("H132").Value = 1843.25
$ws.Range("I132").Value = 1514.8948
$ws.Range("J132").Value = 3922.8333
$ws.Range("K132").Value = 4544.6844
$ws.Range("L132").Value = 11768.4999
$ws.Range("M132").Value = -2014.6844
$ws.Range("N132").Value = -16828.4999

# Row 137
$ws.Range("H137").Value = 2635.853
$ws.Range("I137").Value = 1796.421
$ws.Range("K137").Value = 5389.263
$ws.Range("M137").Value = -2839.263

# Row 138
$ws.Range("H138").Value = 2579.6135
$ws.Range("I138").Value = 929.6667
$ws.Range("J138").Value = 4559.55
$ws.Range("K138").Value = 2789.0001
$ws.Range("L138").Value = 13678.65
$ws.Range("M138").Value = 2350.9999
$ws.Range("N138").Value = -23958.65

# Row 141
$ws.Range("H141").Value = 5207.7
$ws.Range("I141").Value = 5009.75
$ws.Range("K141").Value = 15029.25
$ws.Range("M141").Value = -9849.25

$ws = $wb.Worksheets.Item("ARM")
# Row 9
$ws.Range("H9").Value = 49999
$ws.Range("J9").Value = 49999
$ws.Range("L9").Value = 49999
$ws.Range("N9").Value = -50339

# Row 19
$ws.Range("H19").Value = 24999.5
$ws.Range("I19").Value = 15000
$ws.Range("J19").Value = 34999
$ws.Range("K19").Value = 15000
$ws.Range("L19").Value = 34999
$ws.Range("M19").Value = -14771
$ws.Range("N19").Value = -35457

# Row 20
$ws.Range("H20").Value = 49999
$ws.Range("J20").Value = 49999
$ws.Range("L20").Value = 49999
$ws.Range("N20").Value = -50539

# Row 74
$ws.Range("H74").Value = 80670.39
$ws.Range("I74").Value = 62087.133
$ws.Range("J74").Value = 115514
$ws.Range("K74").Value = 62087.133
$ws.Range("L74").Value = 115514
$ws.Range("M74").Value = -61213.133
$ws.Range("N74").Value = -117262

# Row 77
$ws.Range("H77").Value = 80670.39
$ws.Range("I77").Value = 62087.133
$ws.Range("J77").Value = 115514
$ws.Range("K77").Value = 310435.665
$ws.Range("L77").Value = 577570
$ws.Range("M77").Value = -306067.665
$ws.Range("N77").Value = -586306

# Row 121
$ws.Range("H121").Value = 60038.75
$ws.Range("J121").Value = 60038.75
$ws.Range("L121").Value = 60038.75
$ws.Range("N121").Value = -63532.75

# Row 132
$ws.Range("H132").Value = 8885.228999999999
$ws.Range("I132").Value = 9574.450999999999
$ws.Range("J132").Value = 3543.75
$ws.Range("K132").Value = 28723.353
$ws.Range("L132").Value = 10631.25
$ws.Range("M132").Value = -26193.353
$ws.Range("N132").Value = -15691.25

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 66074
$ws.Range("I86").Value = 1359.4
$ws.Range("K86").Value = 1359.4
$ws.Range("M86").Value = -236.4000000000001

# Row 89
$ws.Range("H89").Value = 66074
$ws.Range("I89").Value = 1359.4
$ws.Range("K89").Value = 6797
$ws.Range("M89").Value = -1181

# Row 105
$ws.Range("H105").Value = 6000
$ws.Range("I105").Value = 6000
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 6000
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -4253
$ws.Range("N105").ClearContents()

# Row 134
$ws.Range("H134").Value = 21890.6
$ws.Range("I134").Value = 25334.842
$ws.Range("J134").Value = 3193.2856
$ws.Range("K134").Value = 76004.526
$ws.Range("L134").Value = 9579.856800000001
$ws.Range("M134").Value = -73469.526
$ws.Range("N134").Value = -14649.8568

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2337.394
$ws.Range("I31").Value = 1575
$ws.Range("K31").Value = 1575
$ws.Range("M31").Value = -1280

# Row 34
$ws.Range("H34").Value = 2337.394
$ws.Range("I34").Value = 1575
$ws.Range("K34").Value = 1575
$ws.Range("M34").Value = -1373

# Row 58
$ws.Range("H58").Value = 7465.4375
$ws.Range("I58").Value = 7024.381
$ws.Range("J58").Value = 8307.454
$ws.Range("K58").Value = 7024.381
$ws.Range("L58").Value = 8307.454
$ws.Range("M58").Value = -6821.381
$ws.Range("N58").Value = -8713.454

# Row 98
$ws.Range("H98").Value = 65811
$ws.Range("J98").Value = 65811
$ws.Range("L98").Value = 65811
$ws.Range("N98").Value = -70303

# Row 107
$ws.Range("H107").Value = 93545.60000000001
$ws.Range("J107").Value = 4763.7144
$ws.Range("L107").Value = 4763.7144
$ws.Range("N107").Value = -8603.714400000001

# Row 132
$ws.Range("H132").Value = 3407
$ws.Range("I132").Value = 3392.7144
$ws.Range("J132").Value = 3507
$ws.Range("K132").Value = 10178.1432
$ws.Range("L132").Value = 10521
$ws.Range("M132").Value = -7648.143199999999
$ws.Range("N132").Value = -15581

# Row 136
$ws.Range("H136").Value = 7465.4375
$ws.Range("I136").Value = 7024.381
$ws.Range("J136").Value = 8307.454
$ws.Range("K136").Value = 21073.143
$ws.Range("L136").Value = 24922.362
$ws.Range("M136").Value = -18523.143
$ws.Range("N136").Value = -30022.362

$ws = $wb.Worksheets.Item("CUL")
# Row 24
$ws.Range("H24").Value = 2000
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 2000
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 6000
$ws.Range("N24").Value = -6460
$ws.Range("M24").ClearContents()

# Row 117
$ws.Range("H117").Value = 982.1429000000001
$ws.Range("I117").Value = 946
$ws.Range("J117").Value = 996.6
$ws.Range("K117").Value = 2838
$ws.Range("L117").Value = 2989.8
$ws.Range("M117").Value = 604
$ws.Range("N117").Value = -9873.799999999999

$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Range("H11").Value = 4140500.2
$ws.Range("I11").Value = 4485125
$ws.Range("J11").Value = 5004
$ws.Range("K11").Value = 4485125
$ws.Range("L11").Value = 5004
$ws.Range("M11").Value = -4484986
$ws.Range("N11").Value = -5282

# Row 102
$ws.Range("H102").Value = 38252.64
$ws.Range("I102").Value = 45284.04
$ws.Range("J102").Value = 10667.923
$ws.Range("K102").Value = 45284.04
$ws.Range("L102").Value = 10667.923
$ws.Range("M102").Value = -43662.04
$ws.Range("N102").Value = -13911.923

$ws = $wb.Worksheets.Item("LTW")
# Row 23
$ws.Range("H23").Value = 7108
$ws.Range("I23").Value = 6001.2
$ws.Range("K23").Value = 6001.2
$ws.Range("M23").Value = -5771.2

# Row 35
$ws.Range("H35").Value = 7433.467
$ws.Range("I35").Value = 1877.1818
$ws.Range("J35").Value = 22713.25
$ws.Range("K35").Value = 1877.1818
$ws.Range("L35").Value = 22713.25
$ws.Range("M35").Value = -1541.1818
$ws.Range("N35").Value = -23385.25

# Row 68
$ws.Range("H68").Value = 2492.6428
$ws.Range("I68").Value = 2081.5454
$ws.Range("K68").Value = 2081.5454
$ws.Range("M68").Value = -1332.5454

# Row 71
$ws.Range("H71").Value = 2492.6428
$ws.Range("I71").Value = 2081.5454
$ws.Range("K71").Value = 10407.727
$ws.Range("M71").Value = -6663.726999999999

# Row 122
$ws.Range("H122").Value = 3595.0312
$ws.Range("I122").Value = 3350.5173
$ws.Range("K122").Value = 10051.5519
$ws.Range("M122").Value = -7601.5519

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 1463.44
$ws.Range("I113").Value = 1352.2858
$ws.Range("J113").Value = 1604.909
$ws.Range("K113").Value = 4056.8574
$ws.Range("L113").Value = 4814.727000000001
$ws.Range("M113").Value = -1886.8574
$ws.Range("N113").Value = -9154.727000000001

# Row 136
$ws.Range("H136").Value = 5334.408
$ws.Range("I136").Value = 5010.206
$ws.Range("K136").Value = 15030.618
$ws.Range("M136").Value = -12480.618
